# Rewrites the large 4th paragraph of the document into 7 paragraphs:
#   Para A .. Para F carry the (edited/new) body text, Para G is the
#   original trailing empty run (kept as-is, just re-homed into its own
#   paragraph).
$d = $word.ActiveDocument

$para = $d.Paragraphs.Item(4)
# Range that excludes the paragraph's own trailing pilcrow so the
# replacement below does not swallow/merge the following paragraph.
$body = $d.Range($para.Range.Start, $para.Range.End - 1)

$t1 = "Fundada em 2019, a Segunda Escola tem como objetivo formar profissionais qualificados para o modelo de mercado atual. Os cursos são organizados de forma rotativa, ou seja, a aula de um dia não necessariamente será a continuidade do dia anterior, então mesmo se o curso já estiver sido iniciado, os alunos podem se matricular na turma e caso o aluno perca determinada aula, ele voltará a vê-la no próximo ciclo, desta forma a entrada e saída de alunos é mais dinâmica do que em escolas tradicionais, tornando necessário um controle de turmas (transferências, trancamentos, e cronogramas) eficiente e integrado."
$t2 = "Atualmente a Segunda Escola tem como meio de gestão o Excel, para cadastrar alunos, controlar turmas, aulas, frequência, notas, professores e mensalidades, tornando trabalhosa e impossível a transferência imediata destas informações, sem esforço manual e demanda de tempo hábil."
$t3 = "A empresa deseja um sistema que auxilie sua gestão administrativa. O sistema deve ter, além de características comuns de um sistema acadêmico comum (controle de: matrícula, turma, frequência, nota, mensalidade, contrato) também características que atendam ao modelo de cursos rotativos, sendo essas as que demandam um maior desafio."
$t4 = "O objetivo do sistema é, além de desafogar o trabalho manual e aumentar a qualidade do atendimento ao aluno, também criar ferramentas que ajudam o gestor a planejar o cronograma recorrentemente e  criar ações de marketing de acordo com projeções da quantidade de alunos nas turmas. O resultado esperado é sempre manter as turmas com alunos, bem como o crescimento da escola."
$t5 = "A empresa deseja um sistema que auxilie em tempo real em sua gestão administrativa, que seja flexível para seguir o modelo cíclico que ela possui em caso de transferências de curso, as informações pessoais, e o histórico do aluno possam ser feitas de forma ágil e simplificada, ainda com níveis de acesso dos usuários, gráficos de vendas, e atividades realizadas nos sistemas (cadastro, relatórios, boletos etc.); aumentando desta forma, a agilidade, transparência e formalização para  a gestão do cliente."
$t6 = "Melhorar texto!"

# Join with carriage returns: Word turns embedded CR into new paragraph
# marks, so this single assignment produces 7 paragraphs in place of 1.
# The trailing CR recreates the final (empty-run) paragraph that used to
# close out the original paragraph.
$body.Text = $t1 + "`r" + $t2 + "`r" + $t3 + "`r" + $t4 + "`r" + $t5 + "`r" + $t6 + "`r"

Write-Output ("ParaCount=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("P" + $i + ": " + $d.Paragraphs.Item($i).Range.Text)
}
